# testdata.xlsx — raw-data update
#
# 1. "address" sheet renamed to "deliveryaddress" and becomes the active tab.
# 2. "productinfo" sheet gains productname/productprice/productdescription
#    header columns (B:D) with matching border styling copied from the
#    existing A2 cell.
# 3. "deliveryaddress" sheet: the mobile number / pincode in B2 & C2 are
#    re-entered as text (they were numbers) and a portrait page setup is
#    applied.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. productinfo: add productname / productprice / productdescription
# ---------------------------------------------------------------------
$wsProduct = $wb.Worksheets.Item("productinfo")

# Reuse A2's existing bordered style for the new B:D header/data cells.
$wsProduct.Range("A2").Copy()
$wsProduct.Range("B1:D2").PasteSpecial(-4122)

$wsProduct.Range("B1").Value = "productname"
$wsProduct.Range("C1").Value = "productprice"
$wsProduct.Range("D1").Value = "productdescription"

$wsProduct.Columns("B:D").AutoFit()

$wsProduct.Range("B3").Select()

# ---------------------------------------------------------------------
# 2. address -> deliveryaddress rename + raw-data fixups
# ---------------------------------------------------------------------
$wsAddress = $wb.Worksheets.Item("address")
$wsAddress.Name = "deliveryaddress"

# Mobile number / pincode become text values (leading apostrophe forces
# text entry, matching the quote-prefixed style in the file).
$wsAddress.Range("B2").Value = "'7338214702"
$wsAddress.Range("C2").Value = "'560016"

$wsAddress.PageSetup.Orientation = 1

# Make this the active sheet/tab, with I2 selected.
$wsAddress.Activate()
$wsAddress.Range("I2").Select()
